$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text (e.g. "62.894.57", "6.06"). Some of
# the new values parse as plain numbers ("6.06", "1.00", ...), so Excel would
# otherwise auto-convert them to numeric cells. Force text entry by flipping
# the cell to Text format for the write, then restore the default "Normal"
# cell style so formatting is left untouched - only the literal text value
# changes, matching the source data.

$ws.Range("D2").Value = '62.894.57'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '3.065.28'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.057.67'
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.06'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.43%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = '3.559.20'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '62.907.27'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '3.062.10'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.65%  '
$ws.Range("E31").Value = '  -6.16%  '
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '486.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.02%  '
$ws.Range("D38").Value = '3.133.75'
$ws.Range("E38").Value = '  +2.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0394'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0794'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.09'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '0.0₃0535'
$ws.Range("E46").Value = '  +8.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.109'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.04%  '
